$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking "Price" column cells to be stored as text,
# matching the source data (e.g. "1.000", "0.000007204") rather than
# letting Excel auto-convert them to numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (name, link, price, volume) cell by cell.
$ws.Range("D2").Value = "30.077.61"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").Value = "1.859.32"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "233.08"
$ws.Range("E5").Value = "  -3.62%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.4649"
$ws.Range("E7").Value = "  -2.95%  "
$ws.Range("D8").Value = "0.2798"
$ws.Range("E8").Value = "  -2.97%  "
$ws.Range("D9").Value = "0.06521"
$ws.Range("E9").Value = "  -3.97%  "
$ws.Range("D10").Value = "19.47"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").Value = "0.07807"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "96.35"
$ws.Range("E12").Value = "  -7.69%  "
$ws.Range("D13").Value = "1.860.06"
$ws.Range("E13").Value = "  -4.19%  "
$ws.Range("D14").Value = "5.109"
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").Value = "0.6630"
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").Value = "280.20"
$ws.Range("E16").Value = "  -4.94%  "
$ws.Range("D17").Value = "30.107.85"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "5.486"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "12.53"
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("D21").Value = "2.100.43"
$ws.Range("E21").Value = "  -4.04%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.000007204"
$ws.Range("E22").Value = "  -5.06%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "6.105"
$ws.Range("E24").Value = "  -4.62%  "
$ws.Range("D25").Value = "9.281"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("D26").Value = "166.42"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").Value = "18.80"
$ws.Range("E27").Value = "  -5.07%  "
$ws.Range("D28").Value = "1.898"
$ws.Range("E28").Value = "  -10.27%  "
$ws.Range("D29").Value = "1.329"
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("D30").Value = "0.09523"
$ws.Range("E30").Value = "  -6.12%  "
$ws.Range("D31").Value = "4.416"
$ws.Range("D32").Value = "1.462"
$ws.Range("E32").Value = "  -4.53%  "
$ws.Range("D33").Value = "4.076"
$ws.Range("E33").Value = "  -6.38%  "
$ws.Range("D34").Value = "0.04623"
$ws.Range("E34").Value = "  -4.32%  "
$ws.Range("D35").Value = "0.6979"
$ws.Range("E35").Value = "  -5.53%  "
$ws.Range("D36").Value = "1.090"
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("D37").Value = "2.698"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").Value = "0.01839"
$ws.Range("E38").Value = "  -5.87%  "
$ws.Range("D39").Value = "6.269"
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("D40").Value = "2.505"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("D41").Value = "72.42"
$ws.Range("E41").Value = "  -5.72%  "
$ws.Range("D42").Value = "0.8524"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.906"
$ws.Range("E43").Value = "  -6.17%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "103.87"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("D46").Value = "0.4133"
$ws.Range("E46").Value = "  -5.25%  "
$ws.Range("D47").Value = "999.56"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.160"
$ws.Range("E48").Value = "  -5.02%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.223"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").Value = "33.88"
$ws.Range("E50").Value = "  -3.51%  "
$ws.Range("E51").Value = "  -6.57%  "
